$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 10:17"

# Swap country names: row 210 was "Islas Malvinas", row 211 was "Groenlandia".
# After the edit, row 210 should show "Groenlandia" and row 211 "Islas Malvinas".
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# Update country statistics rows
# Row 4: Estados Unidos
$ws.Range("B4").Value = 4248492
$ws.Range("C4").Value = 165
$ws.Range("D4").Value = 2028361
$ws.Range("E4").Value = 2071639
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 148492

# Row 6: India
$ws.Range("B6").Value = 1339067
$ws.Range("C6").Value = 2045
$ws.Range("D6").Value = 850295
$ws.Range("E6").Value = 457347
$ws.Range("G6").Value = 19
$ws.Range("H6").Value = 31425

# Row 7: Rusia
$ws.Range("B7").Value = 806720
$ws.Range("C7").Value = 5871
$ws.Range("D7").Value = 597140
$ws.Range("E7").Value = 196388
$ws.Range("G7").Value = 146
$ws.Range("H7").Value = 13192

# Row 21: Alemania
$ws.Range("B21").Value = 205968
$ws.Range("C21").Value = 8
$ws.Range("E21").Value = 6367

# Row 38: Ucrania
$ws.Range("B38").Value = 63929
$ws.Range("C38").Value = 1106
$ws.Range("D38").Value = 35497
$ws.Range("E38").Value = 26842
$ws.Range("G38").Value = 19
$ws.Range("H38").Value = 1590

# Row 53: Armenia
$ws.Range("B53").Value = 36996
$ws.Range("C53").Value = 383
$ws.Range("D53").Value = 26243
$ws.Range("E53").Value = 10053
$ws.Range("G53").Value = 8
$ws.Range("H53").Value = 700

# Row 54: Afganistan
$ws.Range("B54").Value = 36036
$ws.Range("C54").Value = 55
$ws.Range("D54").Value = 24793
$ws.Range("E54").Value = 9997
$ws.Range("G54").Value = 21
$ws.Range("H54").Value = 1246

# Row 56: Kirguistan
$ws.Range("B56").Value = 32124
$ws.Range("C56").Value = 877
$ws.Range("D56").Value = 19203
$ws.Range("E56").Value = 11672
$ws.Range("G56").Value = 38
$ws.Range("H56").Value = 1249

# Row 101: Hungria
$ws.Range("B101").Value = 4424
$ws.Range("C101").Value = 26
$ws.Range("D101").Value = 3324
$ws.Range("E101").Value = 504

# Row 123: Eslovaquia
$ws.Range("B123").Value = 2141
$ws.Range("C123").Value = 23
$ws.Range("E123").Value = 536

# Row 126: Estonia
$ws.Range("B126").Value = 2033
$ws.Range("C126").Value = 5
$ws.Range("E126").Value = 49

# Row 127: Lituania
$ws.Range("B127").Value = 2001
$ws.Range("C127").Value = 15
$ws.Range("E127").Value = 305

# Row 139: Letonia
$ws.Range("B139").Value = 1206
$ws.Range("C139").Value = 1
$ws.Range("E139").Value = 130

# Row 144: Georgia
$ws.Range("B144").Value = 1117
$ws.Range("C144").Value = 13
$ws.Range("D144").Value = 917
$ws.Range("E144").Value = 184

# Row 162: Vietnam
$ws.Range("B162").Value = 416
$ws.Range("C162").Value = 3
$ws.Range("E162").Value = 51

# Row 178: Islas Feroe
$ws.Range("B178").Value = 192
$ws.Range("C178").Value = 1
$ws.Range("E178").Value = 4
